# Rename imu_feature labels from "ss_*" to "sm_*" (reran processing with
# newly cleaned data -> feature prefix changed from ss_ to sm_).

$wb = $excel.ActiveWorkbook

$renames = @{
    "ss_max_60"      = "sm_max_60"
    "ss_max_position" = "sm_max_position"
    "ss_mean_60"     = "sm_mean_60"
    "ss_number_60"   = "sm_number_60"
    "ss_std_60"      = "sm_std_60"
}

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value()
        if ($val -ne $null -and $renames.ContainsKey([string]$val)) {
            $cell.Value = $renames[[string]$val]
        }
    }
}
